# Apply the update described by the commit:
#  - Insert a new item row ("CONCOR 5MG 30 TAB") as the 2nd item, pushing the
#    following items down by one row.
#  - Update the running item numbers (column A) for the shifted rows.
#  - Update the last item's ("UNBLOCKY SOAP") balance/sale-price/transaction
#    figures.
#  - Update the grand total and restore the alternating row heights used by
#    the report generator.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as TEXT (matching the report's habit of
# storing numeric-looking figures as shared strings) without disturbing the
# cell's existing number format / style. It stages the text in a scratch
# cell far outside the used range, then copies only the *value* (not the
# format) into the destination.
function Set-TextValue {
    param($Worksheet, [string]$Address, [string]$TextValue)
    $scratch = $Worksheet.Range("Z100")
    $scratch.NumberFormat = "@"
    $scratch.Value = $TextValue
    $scratch.Copy()
    $Worksheet.Range($Address).PasteSpecial(-4163)
    $scratch.Clear()
}

# ---------------------------------------------------------------------
# 1. Insert a new row right above the current 2nd item (row 8), which
#    shifts GLYCERIN/HIBIOTIC/UNBLOCKY (rows 8-10) down to rows 9-11 and
#    the totals/footer rows (11-12) down to rows 12-13.
# ---------------------------------------------------------------------
$ws.Rows.Item(8).Insert()

# Copy the formatting (styles, borders, the merged-cell look) of the row
# that now sits right below the freshly inserted blank row, so the new row
# looks identical to the other item rows.
$ws.Range($ws.Cells.Item(9,1), $ws.Cells.Item(9,17)).Copy()
$ws.Range($ws.Cells.Item(8,1), $ws.Cells.Item(8,17)).PasteSpecial(-4122)

# Recreate the merged cells for the new row (A:B, C:G, H:K, L:M, N:O).
$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()

# ---------------------------------------------------------------------
# 2. Fill in the values for the new "CONCOR 5MG 30 TAB" row (item #2).
# ---------------------------------------------------------------------
$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "CONCOR 5MG 30 TAB"
$ws.Range("H8").Value = "6:1"
Set-TextValue -Worksheet $ws -Address "L8" -TextValue "1"
Set-TextValue -Worksheet $ws -Address "N8" -TextValue "72.00"
Set-TextValue -Worksheet $ws -Address "P8" -TextValue "72.0000"
$ws.Range("Q8").Value = "1:0"

# ---------------------------------------------------------------------
# 3. Renumber the items that were pushed down (column A values are plain
#    numbers, not formulas, so they don't auto-renumber on row insert).
# ---------------------------------------------------------------------
$ws.Range("A9").Value = 3    # GLYCERIN-ROTEX SOAP
$ws.Range("A10").Value = 4   # HIBIOTIC 1GM 16 TAB
$ws.Range("A11").Value = 5   # UNBLOCKY SOAP

# ---------------------------------------------------------------------
# 4. Update the UNBLOCKY SOAP row (now row 11): balance, sale price and
#    transactions count changed.
# ---------------------------------------------------------------------
$ws.Range("H11").Value = "-1:0"
Set-TextValue -Worksheet $ws -Address "P11" -TextValue "110.0000"
$ws.Range("Q11").Value = "2:0"

# ---------------------------------------------------------------------
# 5. Update the grand total row (now row 12).
# ---------------------------------------------------------------------
$ws.Range("P12").Value = 404.5

# ---------------------------------------------------------------------
# 6. Restore the exact row heights used in the published report (the
#    generator alternates 25.5/24.75 for item rows and keeps 25.5 for the
#    totals row and 16.5 for the footer row).
# ---------------------------------------------------------------------
$ws.Rows.Item(7).RowHeight = 25.5
$ws.Rows.Item(8).RowHeight = 24.75
$ws.Rows.Item(9).RowHeight = 25.5
$ws.Rows.Item(10).RowHeight = 24.75
$ws.Rows.Item(11).RowHeight = 25.5
$ws.Rows.Item(12).RowHeight = 25.5
$ws.Rows.Item(13).RowHeight = 16.5

Write-Output "done"
